{"js": "// Update the \"Updated:\" date line from 2022-05-17 to 2022-07-09\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  if (para.text.indexOf(\"Updated:\") !== -1 && para.text.indexOf(\"2022-05-17\") !== -1) {\n    const searchResults = para.search(\"2022-05-17\", { matchCase: true });\n    searchResults.load(\"items\");\n    await context.sync();\n    for (let j = 0; j < searchResults.items.length; j++) {\n      searchResults.items[j].insertText(\"2022-07-09\", Word.InsertLocation.replace);\n    }\n    await context.sync();\n  }\n}\n", "ps1": "# Update the \"Updated:\" date line from 2022-05-17 to 2022-07-09\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"2022-05-17\", $true, $false, $false, $false, $false, $true, 1, $false, \"2022-07-09\", 2) | Out-Null\n"}
